# Vega Monumental Concepción - Brócoli: add a new weekly entry.
# Two new data rows are inserted at row 265 (pushing the existing rows
# 265..365 down to 267..367), for a new reporting date (serial 44841,
# i.e. 2022-10-07), with a "Primera" and "Segunda" quality record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows, shifting everything at/after row 265 down.
$ws.Rows.Item(265).EntireRow.Insert()
$ws.Rows.Item(266).EntireRow.Insert()

# New row 265: "Primera" quality record for date 44841.
$ws.Range("A265").Value = 11
$ws.Range("B265").Value = "Vega Monumental Concepción"
$ws.Range("C265").Value = "Bíobío"
$ws.Range("D265").Value = 44841
$ws.Range("E265").Value = 8
$ws.Range("F265").Value = 100112023
$ws.Range("G265").Value = "Brócoli"
$ws.Range("H265").Value = "Sin especificar"
$ws.Range("I265").Value = "Primera"
$ws.Range("J265").Value = 2000
$ws.Range("K265").Value = 900
$ws.Range("L265").Value = 1000
$ws.Range("M265").Value = 950
$ws.Range("N265").Value = "$/unidad"
$ws.Range("O265").Value = "Región Metropolitana"
$ws.Range("P265").Value = 950
$ws.Range("Q265").Value = 1
$ws.Range("R265").Value = "Hortaliza"

# New row 266: "Segunda" quality record for the same date.
$ws.Range("A266").Value = 11
$ws.Range("B266").Value = "Vega Monumental Concepción"
$ws.Range("C266").Value = "Bíobío"
$ws.Range("D266").Value = 44841
$ws.Range("E266").Value = 8
$ws.Range("F266").Value = 100112023
$ws.Range("G266").Value = "Brócoli"
$ws.Range("H266").Value = "Sin especificar"
$ws.Range("I266").Value = "Segunda"
$ws.Range("J266").Value = 1000
$ws.Range("K266").Value = 700
$ws.Range("L266").Value = 700
$ws.Range("M266").Value = 700
$ws.Range("N266").Value = "$/unidad"
$ws.Range("O266").Value = "Región Metropolitana"
$ws.Range("P266").Value = 700
$ws.Range("Q266").Value = 1
$ws.Range("R266").Value = "Hortaliza"
